$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Задачи")

# Sheet "Задачи": columns are № | Название | создана | начало | завершено
# The date/time cells in columns C:E share the custom format used throughout
# the sheet (same format already applied to the surrounding cells).
$dateFormat = "dd/mm/yy\ h:mm;@"

# TASK_23 = row 24 ("4.4 - Ход игроком. (LogicGame.doMove)") - completed.
# Row 23 ("4.3 - Покинуть игру. (LogicGame.closeGame)") also records a new
# completion timestamp in column E (завершено).
$ws.Range("E23").Value = 41996.859027777777
$ws.Range("E23").NumberFormat = $dateFormat

# Row 24: record the start (начало, column D) and completion (завершено,
# column E) timestamps for TASK_23 - doMove realized.
$ws.Range("D24").Value = 41997.734027777777
$ws.Range("D24").NumberFormat = $dateFormat

$ws.Range("E24").Value = 41998.57708333333
$ws.Range("E24").NumberFormat = $dateFormat

# Move / leave the active selection on the last edited cell, E24.
$ws.Range("E24").Select()
